$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndices = @(2,3,4,5,7,8,9,10,12,13,14)
$rowNumbers = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25)

$table = @(
    @(2.236963703934521,0.1868048338275798,0.4125214146171032,0.1061728715418604,2.403171259473737,1.950001808201165,1.790789228615473,0.03565943707513242,0.5301895619529802,0.56468599846405,2.277037435095849),
    @(2.152429527502534,0.1663017038063686,0.412753840363365,0.1065940505272485,2.395937815654861,1.953567051358476,1.796697837749022,0.03402759037572167,0.5271813717441916,0.5509545458252774,2.301193869387239),
    @(2.10161077907145,0.1537629123999409,0.4130505900932988,0.1068685861689156,2.392849898295253,1.956644586090903,1.801180497525813,0.03301240256537241,0.5255596765048551,0.5427954884793991,2.316758791818961),
    @(2.081175142173322,0.1486657439627095,0.4132102966172795,0.1069844760972121,2.391931295829508,1.958122009684232,1.803221950636612,0.03259535869621288,0.5249555434970432,0.5395391201019351,2.323285783888581),
    @(2.077798347169789,0.1478201135474535,0.4132391592120328,0.1070039622488577,2.391799265613713,1.958380818668346,1.803573895913395,0.03252590615667472,0.5248586551620917,0.539002543383134,2.324380705833709),
    @(2.101334069292022,0.153694119803589,0.4130525868783579,0.1068701328350738,2.392836134853539,1.956663607149096,1.801207160105648,0.03300679174305898,0.5255512992192735,0.5427512944279016,2.31684607160615),
    @(2.207591235621294,0.1797248442674686,0.4125696096437252,0.1063147930114305,2.400395896860829,1.951046634365355,1.792648989129702,0.03509950894076752,0.5291055984686182,0.5598949643562108,2.285214426840167),
    @(2.424569415025928,0.2311782611702426,0.4128431664678089,0.1053517764760423,2.425990220001097,1.947088092633777,1.782659015553165,0.03909945014675031,0.5378621266166306,0.5956713056003338,2.229002715509409),
    @(2.589244437818763,0.2692454732429326,0.413786771703144,0.1047205135377571,2.451408396584668,1.948493677716272,1.779476985921669,0.04197661094352867,0.5453842042976902,0.6232733762687346,2.19125426450821),
    @(2.6653061770985,0.2866243781544426,0.4143770243399558,0.1044497795810866,2.464418947045374,1.950072754216308,1.778936000888777,0.04327251588355097,0.5490426096769454,0.636116966818733,2.174853240036789),
    @(2.694274131169664,0.2932144442733602,0.4146236590536461,0.1043496136357951,2.469554678019563,1.950806036183764,1.778861792571661,0.04376140780046001,0.5504619482714048,0.6410217925208954,2.168753587976674),
    @(2.688028030339467,0.2917947506919347,0.4145695141543086,0.1043710815394125,2.468439302876476,1.950642089202205,1.778871959212807,0.04365619765381012,0.5501547576519954,0.6399636176526613,2.170062314504143),
    @(2.667686082904197,0.2871663646031095,0.4143968519764769,0.1044414917114511,2.464837275418489,1.95013036860496,1.778927275594313,0.04331277405921696,0.5491586987287462,0.636519663739783,2.174349191015395),
    @(2.655247532751559,0.2843325276759003,0.414294101198891,0.1044849264579506,2.462658159031264,1.949834553553217,1.778978181847151,0.04310217827838159,0.548553007925392,0.6344155106670044,2.176989495993674),
    @(2.584296737337354,0.2681109796641863,0.4137514342316564,0.104738536616158,2.450587320964843,1.948409406577269,1.779530607249036,0.04189166204525208,0.5451498759906173,0.6224397917707236,2.192341640209289),
    @(2.541065052755926,0.2581755519796616,0.4134597382843026,0.1048983211122398,2.443553566347504,1.947775934929922,1.780101894238896,0.04114575218746097,0.5431227286947262,0.6151666003935716,2.201957279086294),
    @(2.516307619511224,0.2524667779208301,0.4133071186716677,0.1049917719150399,2.439644155711704,1.947500012775578,1.780515793121999,0.04071550632802357,0.5419790369413136,0.6110102964470769,2.207560491977929),
    @(2.507943788442503,0.2505348805256347,0.4132580484022839,0.105023678679943,2.438343869743704,1.947421772982921,1.780670575368596,0.04056962262744435,0.5415956286398398,0.6096076903250562,2.209470098397368),
    @(2.545655934774913,0.2592325917000267,0.4134892214514139,0.1048811517279908,2.44428821808205,1.947834215111897,1.780032248459555,0.04122528150122662,0.5433362174422172,0.6159380456850272,2.200926167682486),
    @(2.673656530201129,0.2885255869517493,0.4144469398631117,0.1044207466797515,2.465889602270295,1.950276999206238,1.778907479647373,0.04341369564725994,0.549450343479208,0.6375301178956505,2.173087015704922),
    @(2.758274040879257,0.3077230840256391,0.4152076038776897,0.1041335693860641,2.481225439806991,1.952662343591442,1.778934051383303,0.04483324019292922,0.5536443220547369,0.6518821130011716,2.155539960303617),
    @(2.713024206679506,0.2974721357818169,0.4147893051643763,0.1042855880521527,2.472928712484787,1.951316994653723,1.77885007852592,0.04407657589121072,0.5513878082478811,0.6442002162068334,2.164845845457535),
    @(2.543580094722699,0.258754693977977,0.4134758451368157,0.1048889090579763,2.44395566348183,1.947807591701775,1.78006346914087,0.04118933069990049,0.5432396314534031,0.6155891969224427,2.201392099370917),
    @(2.364948177076542,0.2172134933536825,0.4126386422769741,0.1055988655386138,2.417908327252974,1.947402314590875,1.78463267205143,0.03909945014675031,0.5353020461624567,0.585761687934145,2.243586285244579)
)

for ($r = 0; $r -lt $table.Count; $r++) {
    $rowData = $table[$r]
    $rowNum = $rowNumbers[$r]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $colNum = $colIndices[$c]
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$c]
    }
}
